$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("items")
$tbl = $ws.ListObjects.Item(1)

$ws.Columns.Item(8).Insert()
$tbl.Resize($ws.Range("A1:I25"))

$ws.Range("H1").Value = "productId"
$ws.Range("I1").Value = "link"

$ws.Range("H2").Value = 17
$ws.Range("H3").Value = 9
$ws.Range("H4").Value = 13
$ws.Range("H5").Value = 4
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 3
$ws.Range("H8").Value = 16
$ws.Range("H9").Value = 11
$ws.Range("H10").Value = 3
$ws.Range("H11").Value = 2
$ws.Range("H12").Value = 17
$ws.Range("H13").Value = 8
$ws.Range("H14").Value = 15
$ws.Range("H15").Value = 6
$ws.Range("H16").Value = 12
$ws.Range("H17").Value = 17
$ws.Range("H18").Value = 17
$ws.Range("H19").Value = 10
$ws.Range("H20").Value = 18
$ws.Range("H21").Value = 14
$ws.Range("H22").Value = 12
$ws.Range("H23").Value = 5
$ws.Range("H24").Value = 7
$ws.Range("H25").Value = 8

$ws.Range("I2:I25").Value = "http://www.google.com"
